# ---------------------------------------------------------------------------
# Adds a 4th ("level=4") data row to both the "level 1-3 done" daughter-option
# block (rows 8-14) and the benchmark-diff block (rows 18-24) on the
# "summary" sheet, matching the upstream commit that filled in the
# previously-blank row 11 / row 21 with the newly finished run's numbers,
# then applied a tighter "0.0000" display format to the whole B:G numeric
# block (bold+bordered for the RMSD column, regular+bordered elsewhere).
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("summary")

# --- 1. Fill in the previously-empty "4" rows (11 and 21) ------------------

$ws.Range("B11").Value = 20.0284149341161
$ws.Range("C11").Value = 10.9255241845361
$ws.Range("D11").Value = 5.5318753351538001
$ws.Range("E11").Value = 2.6290960532929
$ws.Range("F11").Value = 1.25806354698319
$ws.Range("G11").Formula = "=SQRT(0.2*SUM(K11:O11))"
$ws.Range("H11").Value = 5798.2030000000004

$ws.Range("B21").Value = 0.105811790323906
$ws.Range("C21").Value = 0.55484458141933402
$ws.Range("D21").Value = 1.3667324235218801
$ws.Range("E21").Value = 2.4156734220249998
$ws.Range("F21").Value = 3.0614458387465699
$ws.Range("G21").Formula = "=SQRT(0.2*SUM(K21:O21))"
$ws.Range("H21").Value = 5798.2030000000004

# --- 2. Re-apply the 4-decimal number format across the whole data block ---
# (B:F on every data row, plus G on the rows that are empty / non-formula;
#  the formula cells in G get the same format applied afterwards together
#  with Bold so they land on their own style, same as before this edit.)

$dataRange = $ws.Range("B8:G14,B18:G24")
$dataRange.NumberFormat = "0.0000"

$rmsdRange = $ws.Range("G8:G11,G18:G21")
$rmsdRange.NumberFormat = "0.0000"
$rmsdRange.Font.Bold = $true

# --- 3. Nudge the affected columns a bit wider now that they carry 4 
#        decimals of data (mirrors Excel's automatic "best fit" after typing)

$ws.Columns.Item(2).ColumnWidth = 8.72
$ws.Columns.Item(3).ColumnWidth = 8.83
$ws.Range("D1:G1").EntireColumn.ColumnWidth = 8.17

# --- 4. Leave the selection where data entry ended -------------------------

$ws.Range("I21").Select()
